# Regenerate merged AHB file:
#  1. Rename the "_old" / "_new" header-column suffixes to "_FV2404" / "_FV2410"
#     (the diff text column in the middle, "diff", is left untouched).
#  2. Turn the used range A1:U58 into a native Excel table ("Table1") so the
#     header row carries the renamed column names and the sheet gets an
#     autoFilter + tableParts reference.
#  3. Freeze the header row (split under row 1, top-left of the scrolling
#     pane anchored at A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" group of headers (columns A..J) to "_FV2404" ---
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

# Column K stays "diff" (index 11) - nothing to do.

# --- Rename the "_new" group of headers (columns L..U) to "_FV2410" ---
$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 11 + 1).Value = $headersFV2410[$i]
}

# --- 2. Convert the used range into a proper table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze panes above row 2 (keep header row visible while scrolling) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
